$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Enter full marks for the two previously-ungraded discussion assignments
$ws.Range("B9").Value = 10
$ws.Range("B10").Value = 10

# Move the active cell selection to O16, matching the saved cursor position
$ws.Range("O16").Select()
